$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_a"
$ws.Range("B2").Value = "scen_a"
$ws.Range("A3").Value = "model_a"
$ws.Range("B3").Value = "scen_c"

$ws.Range("B4").Select()
